$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.126.92'
$ws.Range("E2").Value = '  -1.14%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.857.31'
$ws.Range("E3").Value = '  -2.81%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.18'
$ws.Range("E5").Value = '  -2.31%  '

$ws.Range("E6").Value = '  +0.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4692'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2819'
$ws.Range("E8").Value = '  -1.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06550'
$ws.Range("E9").Value = '  -2.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.19'
$ws.Range("E10").Value = '  +3.69%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07800'
$ws.Range("E11").Value = '  +0.47%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.20'
$ws.Range("E12").Value = '  -6.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.857.55'
$ws.Range("E13").Value = '  -2.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.092'
$ws.Range("E14").Value = '  -1.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6713'
$ws.Range("E15").Value = '  +0.07%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '284.53'
$ws.Range("E16").Value = '  +3.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.158.23'
$ws.Range("E17").Value = '  -1.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.476'
$ws.Range("E19").Value = '  +1.74%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.60'
$ws.Range("E20").Value = '  -0.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.103.32'
$ws.Range("E21").Value = '  -2.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.000007244'
$ws.Range("E22").Value = '  -3.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.145'
$ws.Range("E24").Value = '  -2.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '168.18'
$ws.Range("E25").Value = '  +0.60%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.307'
$ws.Range("E26").Value = '  -0.51%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.07'
$ws.Range("E27").Value = '  -0.94%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.930'
$ws.Range("E28").Value = '  -7.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.341'
$ws.Range("E29").Value = '  -3.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09631'
$ws.Range("E30").Value = '  -3.64%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.407'
$ws.Range("E31").Value = '  -4.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.470'
$ws.Range("E32").Value = '  -2.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.106'
$ws.Range("E33").Value = '  -3.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04683'
$ws.Range("E34").Value = '  -0.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6968'
$ws.Range("E35").Value = '  -4.40%  '

$ws.Range("E36").Value = '  -2.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9994'
$ws.Range("E37").Value = '  +0.18%  '

$ws.Range("E38").Value = '  -0.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01858'
$ws.Range("E39").Value = '  -2.48%  '

$ws.Range("E40").Value = '  -0.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.510'
$ws.Range("E41").Value = '  -3.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.04'
$ws.Range("E42").Value = '  -2.92%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8617'
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.943'
$ws.Range("E44").Value = '  -0.98%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.21'
$ws.Range("E45").Value = '  -2.20%  '

$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9997'
$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("B47").Value = 'TheSandbox'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4161'
$ws.Range("E47").Value = '  -2.44%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.021.00'
$ws.Range("E48").Value = '  +7.10%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.232'
$ws.Range("E49").Value = '  -2.63%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.124'
$ws.Range("E50").Value = '  +4.55%  '

$ws.Range("E51").Value = '  -2.31%  '
